# Apply the "added product like functionality" update to ProjectTimeline.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 currently holds the "27/3/2024" entry (day 23).
# We change its Day text to the new "26/3/2024" entry, and push the
# original "27/3/2024" data down into a brand-new row 27 (day 24) with
# the new hours/description for the product-like feature.

# Preserve the original row26 values before overwriting them
$origDay = $ws.Cells.Item(26, 2).Value2   # "27/3/2024"

# Copy row 26's formatting down into the new row 27 first, so the new row
# matches the rest of the table (same style as used by the other rows)
$ws.Range($ws.Cells.Item(26, 1), $ws.Cells.Item(26, 4)).Copy()
$ws.Range($ws.Cells.Item(27, 1), $ws.Cells.Item(27, 4)).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update row 26: change the Day column to the new date, keep everything else
$ws.Cells.Item(26, 2).Value = "26/3/2024"

# Add new row 27 with the data that used to belong to the "27/3/2024" entry,
# now describing the newly added product-like / retrieval-separation work
$ws.Cells.Item(27, 1).Value = 24
$ws.Cells.Item(27, 2).Value = $origDay
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = "Added product like functionality, separated products retrieval by user"

# Update the selected cell in the sheet view
$ws.Range("C29").Select()

$wb.Save()
